$d = $word.ActiveDocument

# 1. Update the EObjectImpl proxy hash in the first line of the stack trace text
$d.Content.Find.Execute("EObjectImpl@60c44f6f", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "EObjectImpl@1dfdac1f", 2)

# 2. caseQuery line number
$d.Content.Find.Execute("caseQuery(M2DocEvaluator.java:540)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "caseQuery(M2DocEvaluator.java:543)", 2)

# 3. doSwitch(M2DocEvaluator.java:1038) -> 1084 (replace all 3 occurrences)
$d.Content.Find.Execute("doSwitch(M2DocEvaluator.java:1038)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "doSwitch(M2DocEvaluator.java:1084)", 2)

# 4. caseBlock line number
$d.Content.Find.Execute("caseBlock(M2DocEvaluator.java:1254)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "caseBlock(M2DocEvaluator.java:1300)", 2)

# 5. caseDocumentTemplate line number
$d.Content.Find.Execute("caseDocumentTemplate(M2DocEvaluator.java:275)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "caseDocumentTemplate(M2DocEvaluator.java:278)", 2)

# 6. generate(M2DocEvaluator.java:264) -> 267
$d.Content.Find.Execute("generate(M2DocEvaluator.java:264)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "generate(M2DocEvaluator.java:267)", 2)

# 7. M2DocUtils.generate line number
$d.Content.Find.Execute("generate(M2DocUtils.java:712)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "generate(M2DocUtils.java:694)", 2)

# 8. prepareoutputAndGenerate line number
$d.Content.Find.Execute("prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:476)", 2)

# 9. generation line number
$d.Content.Find.Execute("generation(AbstractTemplatesTestSuite.java:369)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "generation(AbstractTemplatesTestSuite.java:385)", 2)
